$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'sliding pads for exercise'
$ws.Cells.Item(2, 1).Value = 'compression capri men'
$ws.Cells.Item(3, 1).Value = 'basketball pants for boys'
$ws.Cells.Item(4, 1).Value = 'padded baseball pants'
$ws.Cells.Item(5, 1).Value = 'cycling capri pants'
$ws.Cells.Item(6, 1).Value = 'sliding pants baseball mens'
$ws.Cells.Item(7, 1).Value = 'below the knee shorts for men'
$ws.Cells.Item(8, 1).Value = 'mens long basketball shorts below knee'
$ws.Cells.Item(9, 1).Value = 'youth small knee pads basketball'
$ws.Cells.Item(10, 1).Value = 'basketball leg pads'
$ws.Cells.Item(11, 1).Value = 'short baseball pants'
$ws.Cells.Item(12, 1).Value = 'basketball bump'
$ws.Cells.Item(13, 1).Value = 'knee pads running'
$ws.Cells.Item(14, 1).Value = 'football compression shorts with pads'
$ws.Cells.Item(15, 1).Value = 'knee high baseball pants mens'
$ws.Cells.Item(16, 1).Value = 'knee padded compression'
$ws.Cells.Item(17, 1).Value = 'athletic pads'
$ws.Cells.Item(18, 1).Value = 'mens protection pads'
$ws.Cells.Item(19, 1).Value = 'youth knee pads wrestling'
$ws.Cells.Item(20, 1).Value = 'men softball pants'
$ws.Cells.Item(21, 1).Value = 'baseball sliding pants'
$ws.Cells.Item(22, 1).Value = 'boys sliding pants'
$ws.Cells.Item(23, 1).Value = 'youth boys basketball pants'
$ws.Cells.Item(24, 1).Value = 'lacrosse sweat pants'
$ws.Cells.Item(25, 1).Value = 'raceface knee pads'
$ws.Cells.Item(26, 1).Value = 'capri compression pants men'
$ws.Cells.Item(27, 1).Value = 'snowboarding compression pants'
$ws.Cells.Item(28, 1).Value = 'men knee pad pants'
$ws.Cells.Item(29, 1).Value = 'youth basketball tights for boys'
$ws.Cells.Item(30, 1).Value = 'knee pads for basketball youth'
$ws.Cells.Item(31, 1).Value = 'mtn bike knee pads'
$ws.Cells.Item(32, 1).Value = 'baseball padded sliding shorts'
$ws.Cells.Item(33, 1).Value = 'leggings with baseballs'
$ws.Cells.Item(34, 1).Value = 'youth knee pad wrestling'
$ws.Cells.Item(35, 1).Value = '28 basketball'
$ws.Cells.Item(36, 1).Value = 'tights with knee'
$ws.Cells.Item(37, 1).Value = 'compression capris'
$ws.Cells.Item(38, 1).Value = 'padded tights for football'
$ws.Cells.Item(39, 1).Value = 'baseball tights for boys'
$ws.Cells.Item(40, 1).Value = 'baseball sliding shorts men'
$ws.Cells.Item(41, 1).Value = 'knee pads for basketball youth boys'
$ws.Cells.Item(42, 1).Value = 'youth leggings boys basketball'
$ws.Cells.Item(43, 1).Value = 'basketballs leggings'
$ws.Cells.Item(44, 1).Value = 'softball items'
$ws.Cells.Item(45, 1).Value = 'football leggings for men'
$ws.Cells.Item(46, 1).Value = 'basketball compression pants youth'
$ws.Cells.Item(47, 1).Value = 'compression pants men black'
$ws.Cells.Item(48, 1).Value = 'mens work pants with knee pads'
$ws.Cells.Item(49, 1).Value = 'youth tights'
$ws.Cells.Item(50, 1).Value = 'mens down pants'
$ws.Cells.Item(51, 1).Value = 'compression basketball pants youth'
$ws.Cells.Item(52, 1).Value = 'basketball compression pants women'
$ws.Cells.Item(53, 1).Value = 'softball slider'
$ws.Cells.Item(54, 1).Value = 'work knee pads under pants'
$ws.Cells.Item(55, 1).Value = 'calf compression pants'
$ws.Cells.Item(56, 1).Value = 'leg pads basketball'
$ws.Cells.Item(57, 1).Value = 'compression pants men football'
$ws.Cells.Item(58, 1).Value = 'softball shorts men'
$ws.Cells.Item(59, 1).Value = 'mens softball shorts'
$ws.Cells.Item(60, 1).Value = 'mens capri shorts below knee'
$ws.Cells.Item(61, 1).Value = 'compression tights youth'
$ws.Cells.Item(62, 1).Value = 'adidas knee pads'
$ws.Cells.Item(63, 1).Value = 'athletic leggings mens'
$ws.Cells.Item(64, 1).Value = 'knee compression shorts'
$ws.Cells.Item(65, 1).Value = 'compression knee pads pair'
$ws.Cells.Item(66, 1).Value = 'indoor volleyball knee pads'
$ws.Cells.Item(67, 1).Value = 'youth basketball pants boys'
$ws.Cells.Item(68, 1).Value = 'softball sliding shorts girls padded'
$ws.Cells.Item(69, 1).Value = 'lacrosse compression shorts padded'
$ws.Cells.Item(70, 1).Value = 'taken leggings'
$ws.Cells.Item(71, 1).Value = 'mens basketball tights'
$ws.Cells.Item(72, 1).Value = 'black compression pants men'
$ws.Cells.Item(73, 1).Value = 'softball sliding shorts'
$ws.Cells.Item(74, 1).Value = 'girls basketball knee pads youth'
$ws.Cells.Item(75, 1).Value = 'kneepads basketball'
$ws.Cells.Item(76, 1).Value = 'knee protector for construction'
$ws.Cells.Item(77, 1).Value = 'knee pad sleeve basketball'
$ws.Cells.Item(78, 1).Value = 'youth basketball leggings'
$ws.Cells.Item(79, 1).Value = 'medium compression pants'
$ws.Cells.Item(80, 1).Value = 'football knee pads for men'
$ws.Cells.Item(81, 1).Value = 'mens compression leggings'
$ws.Cells.Item(82, 1).Value = 'knee pad for yoga'
$ws.Cells.Item(83, 1).Value = 'boys knee pads basketball'
$ws.Cells.Item(84, 1).Value = 'baseball sliding'
$ws.Cells.Item(85, 1).Value = 'knee pads for biking men'
$ws.Cells.Item(86, 1).Value = 'knee pads girls basketball'
$ws.Cells.Item(87, 1).Value = 'black football pants'
$ws.Cells.Item(88, 1).Value = 'lacrosse pants'
$ws.Cells.Item(89, 1).Value = 'team work softball pants'
$ws.Cells.Item(90, 1).Value = 'long basketball shorts for men below knee'
$ws.Cells.Item(91, 1).Value = 'basketball hex pads'
$ws.Cells.Item(92, 1).Value = 'compression shorts padded basketball'
$ws.Cells.Item(93, 1).Value = 'big boys tights'
$ws.Cells.Item(94, 1).Value = 'basketball tights'
$ws.Cells.Item(95, 1).Value = 'hockey knee pads adult'
$ws.Cells.Item(96, 1).Value = 'padded compression pants football'
$ws.Cells.Item(97, 1).Value = 'hockey hip pads adult'
$ws.Cells.Item(98, 1).Value = 'soccer goalkeeper pads'
$ws.Cells.Item(99, 1).Value = 'athletic mens leggings'
$ws.Cells.Item(100, 1).Value = 'knee pads flexible'
